$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: title and link update
$ws.Range("D9").Value = "데이터 사이언스에 관심이 많은 비전공자(디자이너)입니다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/designer-data-science/#utm_source=rss&utm_medium=rss&utm_campaign=designer-data-science"

# Row 26: title update only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 32: title and link update
$ws.Range("D32").Value = "Fine tuning"
$ws.Range("E32").Value = "https://dodonam.tistory.com/347"

# Row 52: title and link update
$ws.Range("D52").Value = "3판 맛보기) 측정단위"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2390&utm_source=rss&utm_medium=rss&utm_campaign=3%25ed%258c%2590-%25eb%25a7%259b%25eb%25b3%25b4%25ea%25b8%25b0-%25ec%25b8%25a1%25ec%25a0%2595%25eb%258b%25a8%25ec%259c%2584"
